$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.235.23'
$ws.Range('D2').Style = $origStyle
$ws.Range('E2').Value = '  +2.75%  '
$origStyle = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.635.81'
$ws.Range('D3').Style = $origStyle
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('E4').Value = '  +0.02%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '567.94'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +5.94%  '
$ws.Range('E6').Value = '  +1.87%  '
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('E8').Value = '  +6.76%  '
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.84'
$ws.Range('D9').Style = $origStyle
$ws.Range('E9').Value = '  -2.68%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.105'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  +4.06%  '
$ws.Range('E11').Value = '  +6.30%  '
$ws.Range('E12').Value = '  +2.27%  '
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.099.47'
$ws.Range('D13').Style = $origStyle
$ws.Range('E13').Value = '  +0.09%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '60.198.82'
$ws.Range('D14').Style = $origStyle
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.71'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  +3.93%  '
$origStyle = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.646.23'
$ws.Range('D16').Style = $origStyle
$ws.Range('E16').Value = '  +0.27%  '
$ws.Range('E17').Value = '  +2.82%  '
$ws.Range('E18').Value = '  +3.72%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '343.20'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  +2.54%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.40'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  +2.30%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.35'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  +1.56%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  +0.00%  '
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.53'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  +0.25%  '
$ws.Range('E24').Value = '  +5.10%  '
$ws.Range('E25').Value = '  +1.78%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.993'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  -0.42%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.30'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  +2.06%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0771'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  +4.32%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('E30').Value = '  +3.80%  '
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.13'
$ws.Range('D31').Style = $origStyle
$ws.Range('E31').Value = '  +5.14%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '156.29'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +3.89%  '
$ws.Range('E33').Value = '  +2.29%  '
$ws.Range('E34').Value = '  +4.84%  '
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.911'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  +7.90%  '
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.912'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  +12.07%  '
$ws.Range('E37').Value = '  +5.38%  '
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '37.42'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  +0.74%  '
$ws.Range('E39').Value = '  +5.66%  '
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '304.20'
$ws.Range('D40').Style = $origStyle
$ws.Range('E40').Value = '  +7.82%  '
$ws.Range('E41').Value = '  +2.31%  '
$ws.Range('E42').Value = '  -0.42%  '
$ws.Range('E43').Value = '  +0.65%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0976'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  +4.16%  '
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0549'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  +3.14%  '
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.35'
$ws.Range('D46').Style = $origStyle
$ws.Range('E46').Value = '  +1.02%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.64'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  -0.47%  '
$ws.Range('E48').Value = '  +4.53%  '
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.52'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  +10.37%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.69'
$ws.Range('D50').Style = $origStyle
$ws.Range('E50').Value = '  +5.35%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.965.13'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  +0.89%  '
